$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The "HORAS" summary block (header + 5 per-person SUMIF rows) currently sits
# at rows 178-183, separated from the data table above (which ends at row
# 174) by three blank rows (175-177). Remove two of those blank rows so the
# block shifts up and only a single blank separator row remains above it.
# ---------------------------------------------------------------------------
$ws.Range("A176:A177").EntireRow.Delete()

# ---------------------------------------------------------------------------
# Append the new list of pending tasks below the summary block, leaving a
# couple of blank (right-aligned styled) rows as a gap, matching the rest of
# the sheet's list formatting.
# ---------------------------------------------------------------------------
$ws.Range("A182").HorizontalAlignment = -4152
$ws.Range("A183").HorizontalAlignment = -4152

$tasks = @(
  "Diseño de los mapas (Colocar los objetos en Tiled)",
  "Arreglar puerta y llave en release",
  "Disparos distintos según arma",
  "Menús en general (nivel en la pantalla seleccionar y meter pantallas finales, créditos)",
  "HUD (Llaves + mini mapa)",
  "Balas",
  "Modelado",
  "Texturizado",
  "Animaciones",
  "Sombras",
  "Shadder Cartoon",
  "Sonido/Música"
)

$r = 184
foreach ($task in $tasks) {
  $ws.Cells.Item($r, 1).Value = $task
  $r = $r + 1
}

# Reorder: "HUD (Llaves + mini mapa)" belongs right after the first item, not
# further down the list - move it up two rows.
$ws.Range("A188").EntireRow.Delete()
$ws.Range("A185").EntireRow.Insert()
$ws.Cells.Item(185, 1).Value = "HUD (Llaves + mini mapa)"

# ---------------------------------------------------------------------------
# Cosmetic view updates: the sheet was scrolled further down and the new
# last cell selected.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 115
$ws.Range("A196").Select()
